$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell values (rows 3,5,6,7 ; columns C,D,E,F) after the edit.
# These correspond to a rotation of the original row values:
#   row3 <- old row5, row5 <- old row7, row6 <- old row3, row7 <- old row6
$updates = @{
    "C3" = "65"; "D3" = "51"; "E3" = "4"; "F3" = "3";
    "C5" = "5";  "D5" = "10"; "E5" = "0"; "F5" = "0";
    "C6" = "72"; "D6" = "53"; "E6" = "6"; "F6" = "2";
    "C7" = "0";  "D7" = "1";  "E7" = "0"; "F7" = "0";
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number format so the numeric-looking string values are
    # written back as text (matching the original t="str" cell type)
    # rather than being auto-converted into numeric cells.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
